$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.948.07"
$ws.Range("D3").Value = "2.919.65"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'593.97"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("D6").Value = "'145.86"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D9").Value = "'6.85"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").Value = "3.400.94"
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").Value = "60.906.29"
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "'6.69"
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("D18").Value = "2.918.46"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("D19").Value = "'430.44"
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("D20").Value = "'13.37"
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("D21").Value = "'0.681"
$ws.Range("E21").Value = "  +1.79%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Value = "'81.56"
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("D24").Value = "'10.94"
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("D26").Value = "'11.94"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "'2.30"
$ws.Range("E28").Value = "  +4.69%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").Value = "'7.04"
$ws.Range("E31").Value = "  -2.89%  "
$ws.Range("D32").Value = "'26.40"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  +0.51%  "
$ws.Range("D34").Value = "0.0₃0850"
$ws.Range("E34").Value = "  +1.99%  "
$ws.Range("E35").Value = "  +0.53%  "
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").Value = "'3.03"
$ws.Range("E37").Value = "  +2.31%  "
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("D39").Value = "'0.123"
$ws.Range("E39").Value = "  -1.92%  "
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("D42").Value = "'40.18"
$ws.Range("E42").Value = "  -3.18%  "
$ws.Range("D43").Value = "'375.94"
$ws.Range("E43").Value = "  +1.55%  "
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("D45").Value = "2.700.54"
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("D46").Value = "'130.84"
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("D48").Value = "'23.95"
$ws.Range("E48").Value = "  -5.69%  "
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("E50").Value = "  -3.77%  "
$ws.Range("E51").Value = "  +2.15%  "
